$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("État de la doc")
$ws.Range("C11").Value = "Owen"
$ws.Range("D13").Select()
